$wb = $excel.ActiveWorkbook

# --- Productdata sheet: AverageDemand for product "1" increases from 40 to 70 ---
$wsProductdata = $wb.Worksheets.Item("Productdata")
# Re-assert the (empty) StandardDevDemands column so the round-trip through the
# COM runtime keeps these cells blank instead of resolving them to shared-string 0.
$wsProductdata.Range("H2:H11").Value = ""
$wsProductdata.Range("G2").Value = 70

# --- ForecastedAverageDemand sheet: give positive demand in the last 3 (uncertainty-free)
#     periods (rows 9-11, buckets 7-9) for product 1 (column B) ---
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

# --- ForcastedStandardDeviation sheet: matching standard deviations for those same cells ---
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
